$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-03 Tuesday" "2024-09-04 Wednesday"

Replace-Text "60×24=1440" "41×68=2788"
Replace-Text "35×73=2555" "24×46=1104"
Replace-Text "33×90=2970" "32×56=1792"
Replace-Text "26×85=2210" "63×65=4095"
Replace-Text "42×77=3234" "41×74=3034"

Replace-Text "60×14=840" "31×15=465"
Replace-Text "96×63=6048" "45×70=3150"
Replace-Text "66×53=3498" "87×83=7221"
Replace-Text "70×37=2590" "50×22=1100"
Replace-Text "20×33=660" "76×89=6764"

Replace-Text "72×70=5040" "47×30=1410"
Replace-Text "17×63=1071" "22×32=704"
Replace-Text "32×15=480" "88×21=1848"
Replace-Text "91×65=5915" "64×26=1664"
Replace-Text "24×95=2280" "19×21=399"

Replace-Text "98×67=6566" "59×56=3304"
Replace-Text "70×30=2100" "97×32=3104"
Replace-Text "48×79=3792" "51×93=4743"
Replace-Text "62×95=5890" "18×16=288"
Replace-Text "11×14=154" "62×91=5642"

Replace-Text "78×20=1560" "67×67=4489"
Replace-Text "49×30=1470" "87×56=4872"
Replace-Text "11×81=891" "67×81=5427"
Replace-Text "55×54=2970" "15×49=735"
Replace-Text "87×58=5046" "73×24=1752"
